$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Frameworks and Tools" bullet: drop the trailing comma after "Express"
#    "Node, JavaScript, Angular, Express, " -> "Node, JavaScript, Angular, Express "
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(9).Range
$p.Find.Execute("Express, ", $true, $false, $false, $false, $false, $true, 1, $false, "Express ", 2)

# ---------------------------------------------------------------------------
# 2) Cloud bullet: add ", Docker" after "AWS, Azure"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(13).Range
$p.Find.Execute("AWS, Azure", $true, $false, $false, $false, $false, $true, 1, $false, "AWS, Azure, Docker", 2)

# ---------------------------------------------------------------------------
# 3) Aureus Tech Systems second bullet: rewrite the eDiscovery paragraph
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(25).Range
$old = "Design and prototype data integration and loading functionality for a cloud-native eDiscovery product in Azure. Prototypes included moving and transforming data from raw zip and CSV files into structured data in Azure SQL Server and document stores in both Blob storage and Azure VMS. Technologies included a range of Azure services."
$new = "Design and prototype integration functionality for a cloud-native eDiscovery product in Azure. Included projects in Azure container server and assorted serverless Azure offerings. "
$p.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ---------------------------------------------------------------------------
# 4) Freelance heading: remove the "Freelance " title and fix the year 2006 -> 2003
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(49).Range
$p.Find.Execute("Freelance ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$p = $d.Paragraphs.Item(49).Range
$p.Find.Execute("2006", $true, $false, $false, $false, $false, $true, 1, $false, "2003", 2)

Write-Host "done"
